# Fruta / hortaliza, semanal
# Insert a new weekly data point at row 57 (pushing the existing rows 57-65
# down to 58-66) for "Mercado Mayorista Lo Valledor de Santiago" - Coco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 57; everything below (57-65) shifts to 58-66.
$ws.Rows("57:57").Insert()

# Populate the new row 57 with the new weekly observation.
$ws.Range("A57").Value = 6
$ws.Range("B57").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C57").Value = "Metropolitana"
$ws.Range("D57").Value = 44711
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100108
$ws.Range("H57").Value = "Tropicales y subtropicales"
$ws.Range("I57").Value = 100108007
$ws.Range("J57").Value = "Coco"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 200
$ws.Range("N57").Value = 20000
$ws.Range("O57").Value = 21000
$ws.Range("P57").Value = 20500
$ws.Range("Q57").Value = "$/malla 20 unidades"
$ws.Range("R57").Value = "Perú"
$ws.Range("S57").Value = 1025
$ws.Range("T57").Value = 20
